$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New column Q: header date Apr-22 (serial 44652), shown as "mmm-yy" ---
# Row 1 already carries a bold+centered row-level style, so the new header
# cell inherits that formatting automatically once a value is written.
$aprDate = Get-Date -Year 2022 -Month 4 -Day 1 -Hour 0 -Minute 0 -Second 0
$ws.Range("Q1").Value = $aprDate
$ws.Range("Q1").NumberFormat = "mmm-yy"

# --- Row 2 (lasgova) updates ---
$ws.Range("N2").Value = 5212
$ws.Range("O2").Value = 5213
$ws.Range("P2").Value = 5224
$ws.Range("Q2").Value = 5260

# --- Row 3 (lalgova) updates ---
$ws.Range("O3").Value = 14073
$ws.Range("P3").Value = 14087
$ws.Range("Q3").Value = 14108

# --- Row 4 (cpgs) updates - commented out deflator pull, new figures pasted
# in with the Haver "Courier New" pull formatting (plain for most cells,
# centered/wrapped for the flagged revision in O4) ---
$ws.Range("N4").Value = 323638
$ws.Range("O4").Value = 324761
$ws.Range("P4").Value = 325304
$ws.Range("Q4").Value = 323956

foreach ($addr in @("N4", "O4", "P4", "Q4")) {
  $cellFont = $ws.Range($addr).Font
  $cellFont.Name = "Courier New"
  $cellFont.Size = 10
  $cellFont.Color = 0
}

$ws.Range("O4").HorizontalAlignment = -4108  # xlCenter
$ws.Range("O4").VerticalAlignment = -4108    # xlCenter
$ws.Range("O4").WrapText = $true

# Match the saved selection from the source workbook
$ws.Range("Q3").Select()

# Page setup was touched by the resave as well (portrait orientation)
$ws.PageSetup.Orientation = 1

